$wb = $excel.ActiveWorkbook

# Sheet "展览": update 想去人数 (F column) for two rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 11863
$wsExhibit.Range("F8").Value = 11786

# Sheet "演出": update 想去人数 (F column) for two rows
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 574
$wsShow.Range("F4").Value = 2

# Sheet "全部类型": update 想去人数 (F column) for four rows (combined view of both sheets above)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 574
$wsAll.Range("F5").Value = 11863
$wsAll.Range("F8").Value = 2
$wsAll.Range("F11").Value = 11786
